$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.36 = 58977.45 pesos`n✅ 58977.45 pesos = 14.33 = 979.57 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsTasas.Range("N10").Value = 69.62
$wsTasas.Range("O10").Value = 4106.01
$wsTasas.Range("N12").Value = 4116.99
$wsTasas.Range("O12").Value = 68.38
